# info page V1.1 minor changes+small easter egg
#
# Mark several existing "week" tasks as done ("y" in column C), and insert
# two new task rows ("experiment with layout" / "simplefy/fix code") under
# the "week 11" section, both also marked done.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark existing tasks as completed ("y") in column C.
$doneRows = @(11, 12, 15, 16, 19, 20, 23, 26)
foreach ($r in $doneRows) {
    $ws.Range("C$r").Value = "y"
}

# Fill in the two previously-blank rows (27-28, a gap between "week 11" at
# row 26 and the existing row 29) with new task entries - this does NOT
# shift any other rows, it just populates the existing gap.
$ws.Range("B27").Value = "experiment with layout"
$ws.Range("C27").Value = "y"

$ws.Range("B28").Value = "simplefy/fix code"
$ws.Range("C28").Value = "y"

# Mark the existing row 29 ("Let Emer know if someone is out of contact")
# as completed too.
$ws.Range("C29").Value = "y"

# Mark the final two checklist items as completed as well.
$ws.Range("C40").Value = "y"
$ws.Range("C41").Value = "y"

# Reflect the final selection from the edit session.
$ws.Range("A37").Select()
